$wb = $excel.ActiveWorkbook

# Rename "Buyers List" sheet to "Buyers"
$buyersSheet = $wb.Worksheets.Item("Buyers List")
$buyersSheet.Name = "Buyers"

# Remove fill formatting from bordered cells (A2:A20, B2:B20) on the Buyers sheet
$buyersSheet.Range("A2:A20").Interior.Pattern = -4142
$buyersSheet.Range("B2:B20").Interior.Pattern = -4142

# Update the active selection on the Buyers sheet
$buyersSheet.Range("A2").Select()
